$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "HK_G_acc_SD"

$values = @(
    87.114337568058076,
    87.114337568058076,
    87.114337568058076,
    87.114337568058076,
    87.114337568058076,
    87.114337568058076,
    86.932849364791281,
    86.932849364791281,
    86.932849364791281,
    86.932849364791281,
    86.932849364791281,
    86.932849364791281,
    86.932849364791281,
    86.932849364791281,
    86.932849364791281,
    86.932849364791281,
    86.932849364791281,
    86.932849364791281,
    87.295825771324871,
    87.295825771324871,
    87.295825771324871,
    87.114337568058076,
    87.114337568058076,
    87.114337568058076,
    88.203266787658805,
    88.021778584392024,
    88.203266787658805,
    86.932849364791281,
    86.932849364791281,
    86.932849364791281,
    88.203266787658805,
    88.384754990925586,
    87.114337568058076,
    87.114337568058076,
    87.114337568058076,
    87.114337568058076,
    87.114337568058076,
    87.114337568058076,
    88.384754990925586,
    86.932849364791281,
    86.932849364791281,
    86.932849364791281,
    88.021778584392024,
    86.932849364791281,
    86.932849364791281,
    86.932849364791281,
    86.932849364791281,
    86.932849364791281
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
